$wb = $excel.ActiveWorkbook

# New data for the two worksheets (rows correspond to sheet row numbers 2..16).
# Column A (Cutoff index) is unchanged; columns B (Cutoff) and C (Reaction_number)
# get the values from what used to be rows 6..20, and the former rows 17..20 are removed.

$sheet1Data = @{
    2  = @(5, 117)
    3  = @(6, 117)
    4  = @(7, 116)
    5  = @(8, 116)
    6  = @(9, 116)
    7  = @(10, 116)
    8  = @(11, 116)
    9  = @(12, 115)
    10 = @(13, 115)
    11 = @(14, 115)
    12 = @(15, 116)
    13 = @(16, 115)
    14 = @(17, 114)
    15 = @(18, 113)
    16 = @(19, 112)
}

$sheet2Data = @{
    2  = @(5, 605)
    3  = @(6, 605)
    4  = @(7, 606)
    5  = @(8, 605)
    6  = @(9, 604)
    7  = @(10, 604)
    8  = @(11, 603)
    9  = @(12, 605)
    10 = @(13, 603)
    11 = @(14, 605)
    12 = @(15, 604)
    13 = @(16, 605)
    14 = @(17, 604)
    15 = @(18, 605)
    16 = @(19, 602)
}

foreach ($sheetInfo in @(
        @{ Name = "NBR"; Data = $sheet1Data },
        @{ Name = "BAR"; Data = $sheet2Data }
    )) {
    $ws = $wb.Worksheets.Item($sheetInfo.Name)

    foreach ($row in $sheetInfo.Data.Keys) {
        $vals = $sheetInfo.Data[$row]
        $ws.Cells.Item($row, 2).Value = $vals[0]
        $ws.Cells.Item($row, 3).Value = $vals[1]
    }

    # Remove the now-obsolete trailing rows (17-20), from the bottom up.
    $ws.Rows.Item(20).EntireRow.Delete()
    $ws.Rows.Item(19).EntireRow.Delete()
    $ws.Rows.Item(18).EntireRow.Delete()
    $ws.Rows.Item(17).EntireRow.Delete()
}
